$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: F17 cuota = 1.44 -> Acierto, profit = cuota - 1
$ws.Range("G17").Value = "Acierto"
$ws.Range("H17").Value = 0.44

# Row 24: F24 cuota = 2 -> Acierto, profit = cuota - 1
$ws.Range("G24").Value = "Acierto"
$ws.Range("H24").Value = 1
